$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in row 52 with the new weekly report entry (7/7/2025)
$ws.Range("D52").Value = Get-Date -Year 2025 -Month 7 -Day 7 -Hour 0 -Minute 0 -Second 0
$ws.Range("E52").Value = 192
$ws.Range("F52").Value = 734
$ws.Range("G52").Value = 0
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 1012
$ws.Range("J52").Value = "Seguire trabajando en fin de semana (rafael)"

# Update the view to reflect scrolled position / selection after edit
$ws.Application.ActiveWindow.ScrollRow = 24
$ws.Range("E54").Select()
